$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Update existing rows 2-7: Runmode column (C) changes from "Y" to "N"
# ---------------------------------------------------------------------
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Value = "N"
}

# ---------------------------------------------------------------------
# 2) Add new rows 8-11 with new search test cases
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "TestCase_B7"
$ws.Range("B8").Value = "To verify that no search results get displayed if search engine doesn't interpret the query"
$ws.Range("C8").Value = "N"
$ws.Range("D8").Value = "SKIP"

$ws.Range("A9").Value = "TestCase_B8"
$ws.Range("B9").Value = "To verify that search maintains state when user navigates back to search results page from record view page"
$ws.Range("C9").Value = "N"
$ws.Range("D9").Value = "SKIP"

$ws.Range("A10").Value = "TestCase_B9"
$ws.Range("B10").Value = "To verify that sorting is retained when user navigates back to search results page from record view page"
$ws.Range("C10").Value = "N"
$ws.Range("D10").Value = "SKIP"

$ws.Range("A11").Value = "TestCase_B10"
$ws.Range("B11").Value = "To verify that filtering is retained when user navigates back to search results page from record view page"
$ws.Range("C11").Value = "Y"
$ws.Range("D11").Value = "SKIP"

# ---------------------------------------------------------------------
# 3) Copy formatting from existing template cells onto the new rows so
#    the new cells pick up the same look (borders / wrap / fill) as the
#    rest of the table.
# ---------------------------------------------------------------------

# D column (no-wrap, bordered) - copy from D2
$ws.Range("D2").Copy()
$ws.Range("D8:D11").PasteSpecial(-4122)

# B8 uses the plain bordered+wrap style (like B2/B5/B6)
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# B9:B11 use the fill-capable bordered+wrap style (like B7)
$ws.Range("B7").Copy()
$ws.Range("B9:B11").PasteSpecial(-4122)

# A8:A11 and C8:C11 use a bordered, non-wrap style that also carries an
# (empty/"no fill") fill flag - start from the plain bordered style (A2)
# then explicitly clear the interior so the fill-apply flag is recorded.
$ws.Range("A2").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C8:C11").PasteSpecial(-4122)

$ws.Range("A8:A11").Interior.Color = 16777215
$ws.Range("A8:A11").Interior.ColorIndex = -4142
$ws.Range("C8:C11").Interior.Color = 16777215
$ws.Range("C8:C11").Interior.ColorIndex = -4142

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) View / selection updates to match the final workbook state
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B16").Select()
